$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The underlying data for each TRI facility (name, location, 1/3/5/10-mile
# neighbor counts) is re-sorted into a new row order. Write the full row
# tuples back out in the new order (header row 1 is untouched).

$rows = @(
    @{ Row = 2;  Facility = "Chemours Chambers Works";                   Location = "Deepwater, NJ";     M1 = 2;  M3 = 3;  M5 = 10; M10 = 32 },
    @{ Row = 3;  Facility = "Chemours El Dorado";                        Location = "El Dorado, AR";      M1 = 2;  M3 = 2;  M5 = 2;  M10 = 12 },
    @{ Row = 4;  Facility = "ARKEMA, INC.";                              Location = "Calvert City, KY";   M1 = 3;  M3 = 11; M5 = 11; M10 = 13 },
    @{ Row = 5;  Facility = "Honeywell International - Geismar Complex"; Location = "Geismar, LA";        M1 = 4;  M3 = 21; M5 = 31; M10 = 36 },
    @{ Row = 6;  Facility = "Chemours - Corpus Christi Plant";           Location = "Gregory, TX";        M1 = 2;  M3 = 4;  M5 = 6;  M10 = 6 },
    @{ Row = 7;  Facility = "Mexichem Fluor Inc.";                       Location = "Saint Gabriel, LA";  M1 = 5;  M3 = 17; M5 = 22; M10 = 36 },
    @{ Row = 8;  Facility = "Iofina Chemical Inc.";                      Location = "Covington, KY";      M1 = 2;  M3 = 2;  M5 = 15; M10 = 44 },
    @{ Row = 9;  Facility = "Chemours Louisville Works";                 Location = "Louisville, KY";     M1 = 12; M3 = 17; M5 = 32; M10 = 66 },
    @{ Row = 10; Facility = "Daikin America Inc.";                       Location = "Decatur, AL";        M1 = 3;  M3 = 16; M5 = 21; M10 = 26 },
    @{ Row = 11; Facility = "Islechem LLC";                              Location = "Grand Island, NY";   M1 = 1;  M3 = 6;  M5 = 11; M10 = 37 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Facility
    $ws.Cells.Item($r.Row, 2).Value = $r.Location
    $ws.Cells.Item($r.Row, 3).Value = $r.M1
    $ws.Cells.Item($r.Row, 4).Value = $r.M3
    $ws.Cells.Item($r.Row, 5).Value = $r.M5
    $ws.Cells.Item($r.Row, 6).Value = $r.M10
}
